$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Locate the "AddressBook" / "Parser" class-diagram box (Rectangle 62) and
# rename it to "RestaurantBook", shrinking the font to fit the longer name -
# matching the rest of the "Changed diagrams to reflect RestaurantBook" edit.
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($shp.HasTextFrame) {
        $tr = $shp.TextFrame.TextRange
        if ($tr.Text -like "*AddressBook*") {
            $para1 = $tr.Paragraphs(1)
            $para1.Text = "RestaurantBook"

            # Resize the whole box's text (both the renamed line and the
            # "Parser" line below it) from 10.5pt down to 9.5pt.
            $tr.Font.Size = 9.5
        }
    }
}
